$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New parking-spot IDs with the trailing "_<GUID>" suffix stripped off.
# Row 2 on the worksheet corresponds to the first entry here, row 144 to the
# last (row 1 is the "ID" header and is left untouched).
$newIds = @(
    "pbg-s-eladen-1-01",
    "pbg-s-eladen-1-02",
    "pbg-s-eladen-1-03",
    "pbg-s-eladen-1-04",
    "pbg-s-eladen-1-05",
    "pbg-s-eladen-1-06",
    "pbg-s-eladen-1-07",
    "pbg-s-eladen-1-08",
    "pbg-s-eladen-1-09",
    "pbg-s-eladen-1-10",
    "pbg-s-eladen-1-11",
    "pbg-s-eladen-1-12",
    "pbg-s-eladen-1-13",
    "pbg-s-eladen-1-14",
    "pbg-s-eladen-1-15",
    "pbg-s-eladen-1-16",
    "pbg-s-eladen-1-17",
    "pbg-s-eladen-1-18",
    "pbg-s-eladen-1-19",
    "pbg-s-eladen-1-20",
    "pbg-f-eladen-0-01",
    "pbg-f-eladen-0-02",
    "pbg-f-eladen-0-03",
    "pbg-f-eladen-0-04",
    "pbg-f-eladen-0-05",
    "pbg-f-eladen-0-06",
    "pbg-f-eladen-0-07",
    "pbg-f-eladen-0-08",
    "pbg-f-eladen-0-09",
    "pbg-f-eladen-0-10",
    "pbg-f-eladen-0-11",
    "pbg-f-eladen-0-12",
    "pbg-f-eladen-0-13",
    "pbg-f-eladen-0-14",
    "pbg-f-eladen-0-15",
    "pbg-f-eladen-0-16",
    "pbg-f-eladen-0-17",
    "pbg-f-eladen-0-18",
    "pbg-f-eladen-0-19",
    "pbg-f-eladen-0-20",
    "pbg-c01-19296d",
    "pbg-c02-1921b0",
    "pbg-c03-192985",
    "pbg-c04-192995",
    "pbg-c05-1921cd",
    "pbg-c06-1921dd",
    "pbg-c07-192968",
    "pbg-c08-1921da",
    "pbg-c09-192980",
    "pbg-c10-1921bc",
    "pbg-c11-192996",
    "pbg-c12-19298a",
    "pbg-c13-1921bb",
    "pbg-c14-1921d4",
    "pbg-c15-192987",
    "pbg-c16-19297d",
    "pbg-r01-192851",
    "pbg-r02-19284f",
    "pbg-r03-192839",
    "pbg-r04-192835",
    "pbg-r05-192837",
    "pbg-r06-192838",
    "pbg-r07-192832",
    "pbg-r08-192834",
    "pbg-r09-192833",
    "pbg-r10-192836",
    "pbg-r11-192830",
    "pbg-r12-192831",
    "pbg-r13-19283c",
    "pbg-r14-19283d",
    "pbg-r15-19283b",
    "pbg-r16-19283a",
    "pbg-r17-192827",
    "pbg-r18-19283e",
    "pbg-r19-19283f",
    "pbg-r20-192840",
    "pbg-r21-19282f",
    "pbg-r22-19282e",
    "pbg-r23-192841",
    "pbg-r24-19282c",
    "pbg-r25-192842",
    "pbg-r26-19282b",
    "pbg-r27-19282d",
    "pbg-r28-19282a",
    "pbg-r29-192829",
    "pbg-r30-192825",
    "pbg-r31-192826",
    "pbg-r32-192828",
    "pbg-b-elade-o-1-006",
    "pbg-b-elade-o-1-007",
    "pbg-b-elade-o-1-008",
    "pbg-b-elade-o-1-009",
    "pbg-b-elade-o-1-010",
    "pbg-b-elade-o-1-011",
    "pbg-b-elade-o-1-012",
    "pbg-b-elade-o-1-013",
    "pbg-b-elade-o-1-014",
    "pbg-b-elade-o-1-015",
    "pbg-b-elade-o-1-016",
    "pbg-b-elade-o-1-017",
    "pbg-b-elade-o-1-018",
    "pbg-b-elade-o-1-019",
    "pbg-b-elade-o-1-020",
    "pbg-b-elade-o-1-021",
    "pbg-b-elade-w-1-054",
    "pbg-b-elade-w-1-055",
    "pbg-b-elade-w-1-056",
    "pbg-b-elade-w-1-057",
    "pbg-b-elade-w-1-058",
    "pbg-b-elade-w-1-059",
    "pbg-b-elade-w-1-060",
    "pbg-b-elade-w-1-061",
    "pbg-b-elade-w-1-062",
    "pbg-b-elade-w-1-063",
    "pbg-b-elade-w-1-064",
    "pbg-b-elade-w-1-065",
    "pbg-b-elade-w-1-066",
    "pbg-b-elade-w-1-067",
    "pbg-b-elade-w-1-068",
    "pbg-b-elade-w-1-069",
    "pbg-b-familie-o-1-001",
    "pbg-b-familie-o-1-002",
    "pbg-b-familie-o-2-001",
    "pbg-b-familie-o-2-002",
    "pbg-b-familie-o-3-001",
    "pbg-b-familie-o-3-002",
    "pbg-b-familie-o-4-001",
    "pbg-b-familie-o-4-002",
    "pbg-b-familie-w-1-003",
    "pbg-b-familie-w-1-004",
    "pbg-b-familie-w-2-003",
    "pbg-b-familie-w-2-004",
    "pbg-b-familie-w-3-003",
    "pbg-b-familie-w-3-004",
    "pbg-b-familie-w-4-003",
    "pbg-b-familie-w-4-004",
    "pbg-b-handicap-o-1-001",
    "pbg-b-handicap-o-2-001",
    "pbg-b-handicap-o-3-001",
    "pbg-b-handicap-o-4-001",
    "pbg-b-handicap-w-2-002",
    "pbg-b-handicap-w-3-002",
    "pbg-b-handicap-w-4-002"
)

for ($i = 0; $i -lt $newIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newIds[$i]
}

# Restore the view state recorded in the saved workbook: scrolled so row 132
# is at the top, with A141 as the active selection.
$win = $excel.ActiveWindow
$ws.Range("A141").Select()
$win.ScrollRow = 132
$win.ScrollColumn = 1
